$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data cells (B/C/D/E for every row) are stored as plain text
# in the workbook, even numeric-looking prices like "0.996". Excel's COM
# layer auto-coerces a clean numeric string assigned via .Value into a
# Number cell, so B/C/D are temporarily forced to Text format ("@") while
# the new values are written, then the temporary formatting is cleared so
# the cells end up with no explicit style (matching the source file) while
# keeping the stored type as text. Column E values are percentages like
# "  +0.44%  " which never parse as numbers, so they need no such handling.
$ws.Range("B2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "57.063.40"
$ws.Range("E2").Value = "  +0.44%  "

# Row 3
$ws.Range("D3").Value = "2.401.68"
$ws.Range("E3").Value = "  -3.60%  "

# Row 4
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").Value = "486.55"
$ws.Range("E5").Value = "  -0.95%  "

# Row 6
$ws.Range("D6").Value = "154.77"
$ws.Range("E6").Value = "  +1.73%  "

# Row 7
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  +17.99%  "

# Row 9
$ws.Range("D9").Value = "2.420.10"
$ws.Range("E9").Value = "  -3.33%  "

# Row 10
$ws.Range("D10").Value = "6.30"
$ws.Range("E10").Value = "  +10.06%  "

# Row 11
$ws.Range("D11").Value = "0.0995"
$ws.Range("E11").Value = "  +0.89%  "

# Row 12
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +0.27%  "

# Row 13
$ws.Range("E13").Value = "  +1.27%  "

# Row 14
$ws.Range("D14").Value = "2.820.71"
$ws.Range("E14").Value = "  -3.59%  "

# Row 15
$ws.Range("D15").Value = "57.015.72"
$ws.Range("E15").Value = "  +0.53%  "

# Row 16
$ws.Range("E16").Value = "  -2.67%  "

# Row 17
$ws.Range("E17").Value = "  -2.17%  "

# Row 18
$ws.Range("D18").Value = "2.412.42"
$ws.Range("E18").Value = "  -3.42%  "

# Row 19
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  +3.77%  "

# Row 20
$ws.Range("D20").Value = "324.85"
$ws.Range("E20").Value = "  +1.11%  "

# Row 21
$ws.Range("D21").Value = "9.91"
$ws.Range("E21").Value = "  -3.94%  "

# Row 22
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("E23").Value = "  +1.06%  "

# Row 24
$ws.Range("D24").Value = "58.25"
$ws.Range("E24").Value = "  -0.31%  "

# Row 25
$ws.Range("E25").Value = "  -1.56%  "

# Row 26
$ws.Range("E26").Value = "  -0.53%  "

# Row 27
$ws.Range("E27").Value = "  -0.72%  "

# Row 28
$ws.Range("D28").Value = "2.513.22"
$ws.Range("E28").Value = "  -3.16%  "

# Row 29
$ws.Range("D29").Value = "7.24"
$ws.Range("E29").Value = "  -4.42%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0779"
$ws.Range("E30").Value = "  -2.85%  "

# Row 31
$ws.Range("E31").Value = "  +0.03%  "

# Row 32
$ws.Range("D32").Value = "149.98"
$ws.Range("E32").Value = "  -1.22%  "

# Row 33
$ws.Range("D33").Value = "18.54"
$ws.Range("E33").Value = "  +1.19%  "

# Row 34
$ws.Range("E34").Value = "  +0.12%  "

# Row 35
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("D36").Value = "1.16"
$ws.Range("E36").Value = "  -0.20%  "

# Row 37
$ws.Range("D37").Value = "3.73"
$ws.Range("E37").Value = "  -1.06%  "

# Row 38
$ws.Range("D38").Value = "0.839"
$ws.Range("E38").Value = "  -3.45%  "

# Row 39
$ws.Range("D39").Value = "34.10"
$ws.Range("E39").Value = "  -0.26%  "

# Row 40
$ws.Range("D40").Value = "0.101"
$ws.Range("E40").Value = "  +8.24%  "

# Row 41
$ws.Range("E41").Value = "  -0.10%  "

# Row 42
$ws.Range("D42").Value = "1.37"
$ws.Range("E42").Value = "  -1.84%  "

# Row 43
$ws.Range("E43").Value = "  -0.27%  "

# Row 44
$ws.Range("D44").Value = "0.595"
$ws.Range("E44").Value = "  -3.22%  "

# Row 45
$ws.Range("D45").Value = "268.65"
$ws.Range("E45").Value = "  +0.67%  "

# Row 46
$ws.Range("D46").Value = "0.0529"
$ws.Range("E46").Value = "  -6.07%  "

# Row 47
$ws.Range("D47").Value = "10.22"
$ws.Range("E47").Value = "  +0.12%  "

# Row 48
$ws.Range("E48").Value = "  -0.09%  "

# Row 49
$ws.Range("D49").Value = "4.56"
$ws.Range("E49").Value = "  -5.10%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "17.41"
$ws.Range("E50").Value = "  -2.50%  "

# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.867.02"
$ws.Range("E51").Value = "  -1.48%  "

$ws.Range("B2:D51").ClearFormats()